$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Signature placeholders for the approval stamp cells (F/G/H/I, rows 3 and 5),
# added column by column so the shared-string table is populated in the same
# order as the authored workbook (DAM, TIM, BU, DEA - each with SIG1/SIG2).
$ws.Range("F3").Value = "{{DAM_SIG1}}"
$ws.Range("F5").Value = "{{DAM_SIG2}}"
$ws.Range("G3").Value = "{{TIM_SIG1}}"
$ws.Range("G5").Value = "{{TIM_SIG2}}"
$ws.Range("H3").Value = "{{BU_SIG1}}"
$ws.Range("H5").Value = "{{BU_SIG2}}"
$ws.Range("I3").Value = "{{DEA_SIG1}}"
$ws.Range("I5").Value = "{{DEA_SIG2}}"

# Update the selection to match the authored state
$ws.Range("A6:I8").Select()
